$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting the old
# "Late"/"heading"/"Outstanding" columns (N,O,P) one to the right (O,P,Q).
$ws.Columns("N").Insert()

# The newly inserted column picks up the width of the column to its left
# (column M) without the "best fit" auto-size flag, matching how Excel
# formats a manually inserted column.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and select cell R8,
# mirroring the saved selection state in the workbook.
$ws.Activate()
$ws.Range("R8").Select()
